$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet so it lands at the end
# of the tab order (Tabelle1, Sheet2, Sheet3).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Populate the new sheet's data (A1:D2 used range; column C intentionally
# left blank so there is no "Col3").
$ws3.Range("A1").Value = "Col1"
$ws3.Range("B1").Value = "Col2"
$ws3.Range("D1").Value = "Col4"
$ws3.Range("A2").Value = "text2"
$ws3.Range("B2").Value = "text1"
$ws3.Range("D2").Value = "text3"

# Make the new sheet the active/selected tab.
$ws3.Activate()
